$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '54.156.95'
$ws.Range("E2").Value = '  -11.06%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.278.19'
$ws.Range("E3").Value = '  -21.74%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '451.37'
$ws.Range("E5").Value = '  -15.03%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '129.09'
$ws.Range("E6").Value = '  -11.36%  '
$ws.Range("E7").Value = '  -0.15%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.473'
$ws.Range("E8").Value = '  -15.10%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.274.32'
$ws.Range("E9").Value = '  -22.04%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '5.44'
$ws.Range("E10").Value = '  -9.82%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0918'
$ws.Range("E11").Value = '  -15.75%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.310'
$ws.Range("E12").Value = '  -15.35%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.122'
$ws.Range("E13").Value = '  -2.75%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.659.77'
$ws.Range("E14").Value = '  -22.21%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '54.187.13'
$ws.Range("E15").Value = '  -10.90%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '18.91'
$ws.Range("E16").Value = '  -17.42%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000120'
$ws.Range("E17").Value = '  -15.91%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.293.01'
$ws.Range("E18").Value = '  -21.13%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.07'
$ws.Range("E19").Value = '  -19.76%  '
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '304.63'
$ws.Range("E20").Value = '  -16.38%  '
$ws.Range("B21").Value = 'Chainlink'
$ws.Range("C21").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.55'
$ws.Range("E21").Value = '  -18.63%  '
$ws.Range("E22").Value = '  -0.22%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.62'
$ws.Range("E23").Value = '  -1.67%  '
$ws.Range("E24").Value = '  -19.71%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '55.76'
$ws.Range("E25").Value = '  -13.91%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.983'
$ws.Range("E26").Value = '  -1.69%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.160'
$ws.Range("E27").Value = '  -12.99%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.374'
$ws.Range("E28").Value = '  -18.16%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.84'
$ws.Range("E29").Value = '  -13.11%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.995'
$ws.Range("E30").Value = '  -0.37%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0₃0713'
$ws.Range("E31").Value = '  -18.35%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '144.36'
$ws.Range("E32").Value = '  -4.50%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '16.96'
$ws.Range("E33").Value = '  -14.32%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.37'
$ws.Range("E34").Value = '  -19.24%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.76'
$ws.Range("E35").Value = '  -15.14%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.66'
$ws.Range("E36").Value = '  -17.28%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.843'
$ws.Range("E37").Value = '  -16.80%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.00'
$ws.Range("E38").Value = '  -17.46%  '
$ws.Range("B39").Value = 'FirstDigitalUSD'
$ws.Range("C39").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.993'
$ws.Range("E39").Value = '  -0.32%  '
$ws.Range("B40").Value = 'OKB'
$ws.Range("C40").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '33.10'
$ws.Range("E40").Value = '  -12.45%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '10.30'
$ws.Range("E41").Value = '  -0.40%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.25'
$ws.Range("E42").Value = '  -16.73%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.19'
$ws.Range("E43").Value = '  -15.15%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.924.96'
$ws.Range("E44").Value = '  -16.11%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0499'
$ws.Range("E45").Value = '  -14.69%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0206'
$ws.Range("E46").Value = '  -13.63%  '
$ws.Range("B47").Value = 'Mantle'
$ws.Range("C47").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.509'
$ws.Range("E47").Value = '  -21.63%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0811'
$ws.Range("E48").Value = '  -12.85%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '16.37'
$ws.Range("E49").Value = '  -21.60%  '
$ws.Range("B50").Value = 'RenderToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.08'
$ws.Range("E50").Value = '  -19.36%  '
$ws.Range("B51").Value = 'ZEEBU'
$ws.Range("C51").Value = 'https://coinranking.com/coin/B5-YKN_zB+zeebu-zbu'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.68'
$ws.Range("E51").Value = '  -3.03%  '
